$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 20.598495
$ws.Range("H2").Value = 61.795485
$ws.Range("I2").Value = 0.187290770808477
$ws.Range("J2").Value = 0.187290770808477
$ws.Range("M2").Value = 4.265285666666666
$ws.Range("N2").Value = 12.795857
$ws.Range("O2").Value = 0.3060251613083444
$ws.Range("P2").Value = 0.3060251613083444
$ws.Range("Q2").Value = 87.858465478405
$ws.Range("R2").Value = 790.726189305645
$ws.Range("S2").Value = 0.05731568834822835
$ws.Range("T2").Value = 0.05731568834822835

# Row 3
$ws.Range("G3").Value = 20.598495
$ws.Range("H3").Value = 61.795485
$ws.Range("I3").Value = 0.187290770808477
$ws.Range("J3").Value = 0.187290770808477
$ws.Range("O3").Value = 0.5138779274087578
$ws.Range("P3").Value = 0.5138779274087578
$ws.Range("Q3").Value = 147.532072043795
$ws.Range("R3").Value = 1327.788648394155
$ws.Range("S3").Value = 0.09624459312584886
$ws.Range("T3").Value = 0.09624459312584886

# Row 4
$ws.Range("G4").Value = 20.598495
$ws.Range("H4").Value = 61.795485
$ws.Range("I4").Value = 0.187290770808477
$ws.Range("J4").Value = 0.187290770808477
$ws.Range("M4").Value = 2.510136
$ws.Range("N4").Value = 7.530408
$ws.Range("O4").Value = 0.1800969112828978
$ws.Range("P4").Value = 0.1800969112828978
$ws.Range("Q4").Value = 51.70502384531999
$ws.Range("R4").Value = 465.3452146078799
$ws.Range("S4").Value = 0.03373048933439984
$ws.Range("T4").Value = 0.03373048933439984

# Row 5
$ws.Range("G5").Value = 60.20577233333334
$ws.Range("I5").Value = 0.5474179306512287
$ws.Range("J5").Value = 0.5474179306512288
$ws.Range("M5").Value = 4.265285666666666
$ws.Range("N5").Value = 12.795857
$ws.Range("O5").Value = 0.3060251613083444
$ws.Range("P5").Value = 0.3060251613083444
$ws.Range("Q5").Value = 256.7948177839632
$ws.Range("R5").Value = 2311.153360055669
$ws.Range("S5").Value = 0.1675236605306223
$ws.Range("T5").Value = 0.1675236605306224

# Row 6
$ws.Range("G6").Value = 60.20577233333334
$ws.Range("I6").Value = 0.5474179306512287
$ws.Range("J6").Value = 0.5474179306512288
$ws.Range("O6").Value = 0.5138779274087578
$ws.Range("P6").Value = 0.5138779274087578
$ws.Range("Q6").Value = 431.2102579015435
$ws.Range("S6").Value = 0.2813059916294445
$ws.Range("T6").Value = 0.2813059916294446

# Row 7
$ws.Range("G7").Value = 60.20577233333334
$ws.Range("I7").Value = 0.5474179306512287
$ws.Range("J7").Value = 0.5474179306512288
$ws.Range("M7").Value = 2.510136
$ws.Range("N7").Value = 7.530408
$ws.Range("O7").Value = 0.1800969112828978
$ws.Range("P7").Value = 0.1800969112828978
$ws.Range("Q7").Value = 151.124676541704
$ws.Range("R7").Value = 1360.122088875336
$ws.Range("S7").Value = 0.09858827849116183
$ws.Range("T7").Value = 0.09858827849116186

# Row 8
$ws.Range("G8").Value = 29.17709966666666
$ws.Range("H8").Value = 87.53129899999999
$ws.Range("I8").Value = 0.2652912985402942
$ws.Range("J8").Value = 0.2652912985402942
$ws.Range("M8").Value = 4.265285666666666
$ws.Range("N8").Value = 12.795857
$ws.Range("O8").Value = 0.3060251613083444
$ws.Range("P8").Value = 0.3060251613083444
$ws.Range("Q8").Value = 124.4486650031381
$ws.Range("R8").Value = 1120.037985028243
$ws.Range("S8").Value = 0.08118581242949369
$ws.Range("T8").Value = 0.08118581242949369

# Row 9
$ws.Range("G9").Value = 29.17709966666666
$ws.Range("H9").Value = 87.53129899999999
$ws.Range("I9").Value = 0.2652912985402942
$ws.Range("J9").Value = 0.2652912985402942
$ws.Range("O9").Value = 0.5138779274087578
$ws.Range("P9").Value = 0.5138779274087578
$ws.Range("Q9").Value = 208.9743920636752
$ws.Range("R9").Value = 1880.769528573077
$ws.Range("S9").Value = 0.1363273426534644
$ws.Range("T9").Value = 0.1363273426534644

# Row 10
$ws.Range("G10").Value = 29.17709966666666
$ws.Range("H10").Value = 87.53129899999999
$ws.Range("I10").Value = 0.2652912985402942
$ws.Range("J10").Value = 0.2652912985402942
$ws.Range("M10").Value = 2.510136
$ws.Range("N10").Value = 7.530408
$ws.Range("O10").Value = 0.1800969112828978
$ws.Range("P10").Value = 0.1800969112828978
$ws.Range("Q10").Value = 73.23848824888799
$ws.Range("R10").Value = 659.1463942399919
$ws.Range("S10").Value = 0.04777814345733612
$ws.Range("T10").Value = 0.04777814345733612
